# Dimensionen.xlsx — "Added the positions to the numbers and letters"
#
# - A9 ("Letters") becomes "Direction"
# - The "Bezeichner" block's Scaled Size/x values change (G16, I16),
#   which ripples into the cached result of the G17 formula (=I16*G16)
#   automatically on recalculation.
# - Selection moves to H19 (with the view scrolled so C7 is the
#   top-left visible cell — window scroll position is transient host
#   UI state and outside of what gets written back into the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Letters" -> "Direction"
$ws.Range("A9").Value = "Direction"

# Bezeichner block: Scaled Size 56 -> 43, Scaled x 36 -> 27
$ws.Range("G16").Value = 43
$ws.Range("I16").Value = 27

# Recalculate so the cached formula result in G17 (=I16*G16) updates to 1161
$excel.Calculate()

# Update the view: scroll so column C / row 7 is the top-left corner,
# and select H19
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("H19").Select()
